$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 392.77777  # ALC!H28
$ws.Cells.Item(28, 9).Value = 418.2143  # ALC!I28
$ws.Cells.Item(28, 11).Value = 418.2143  # ALC!K28
$ws.Cells.Item(28, 13).Value = 66.78570000000002  # ALC!M28

$ws.Cells.Item(53, 8).Value = 1286.75  # ALC!H53
$ws.Cells.Item(53, 10).Value = 542.44446  # ALC!J53
$ws.Cells.Item(53, 12).Value = 542.44446  # ALC!L53
$ws.Cells.Item(53, 14).Value = -1816.44446  # ALC!N53

$ws.Cells.Item(62, 8).Value = 31252634  # ALC!H62
$ws.Cells.Item(62, 9).Value = 35716772  # ALC!I62
$ws.Cells.Item(62, 10).Value = 3687  # ALC!J62
$ws.Cells.Item(62, 11).Value = 35716772  # ALC!K62
$ws.Cells.Item(62, 12).Value = 3687  # ALC!L62
$ws.Cells.Item(62, 13).Value = -35716148  # ALC!M62
$ws.Cells.Item(62, 14).Value = -4935  # ALC!N62

$ws.Cells.Item(65, 8).Value = 31252634  # ALC!H65
$ws.Cells.Item(65, 9).Value = 35716772  # ALC!I65
$ws.Cells.Item(65, 10).Value = 3687  # ALC!J65
$ws.Cells.Item(65, 11).Value = 178583860  # ALC!K65
$ws.Cells.Item(65, 12).Value = 18435  # ALC!L65
$ws.Cells.Item(65, 13).Value = -178580740  # ALC!M65
$ws.Cells.Item(65, 14).Value = -24675  # ALC!N65

$ws.Cells.Item(88, 8).Value = 1746.8572  # ALC!H88
$ws.Cells.Item(88, 9).Value = 1611.8572  # ALC!I88
$ws.Cells.Item(88, 10).Value = 1881.8572  # ALC!J88
$ws.Cells.Item(88, 11).Value = 1611.8572  # ALC!K88
$ws.Cells.Item(88, 12).Value = 1881.8572  # ALC!L88
$ws.Cells.Item(88, 13).Value = -1205.8572  # ALC!M88
$ws.Cells.Item(88, 14).Value = -2693.8572  # ALC!N88

$ws.Cells.Item(91, 8).Value = 1746.8572  # ALC!H91
$ws.Cells.Item(91, 9).Value = 1611.8572  # ALC!I91
$ws.Cells.Item(91, 10).Value = 1881.8572  # ALC!J91
$ws.Cells.Item(91, 11).Value = 1611.8572  # ALC!K91
$ws.Cells.Item(91, 12).Value = 1881.8572  # ALC!L91
$ws.Cells.Item(91, 13).Value = -207.8571999999999  # ALC!M91
$ws.Cells.Item(91, 14).Value = -4689.8572  # ALC!N91

$ws.Cells.Item(92, 8).Value = 1516.4231  # ALC!H92
$ws.Cells.Item(92, 9).Value = 575.7368  # ALC!I92
$ws.Cells.Item(92, 10).Value = 4069.7144  # ALC!J92
$ws.Cells.Item(92, 11).Value = 575.7368  # ALC!K92
$ws.Cells.Item(92, 12).Value = 4069.7144  # ALC!L92
$ws.Cells.Item(92, 13).Value = 672.2632  # ALC!M92
$ws.Cells.Item(92, 14).Value = -6565.7144  # ALC!N92

$ws.Cells.Item(94, 8).Value = 999.3333  # ALC!H94
$ws.Cells.Item(94, 9).Value = 999.3333  # ALC!I94
$ws.Cells.Item(94, 11).Value = 999.3333  # ALC!K94
$ws.Cells.Item(94, 13).Value = -548.3333  # ALC!M94

$ws.Cells.Item(96, 8).Value = 796.55554  # ALC!H96
$ws.Cells.Item(96, 9).Value = 844.25  # ALC!I96
$ws.Cells.Item(96, 10).Value = 758.4  # ALC!J96
$ws.Cells.Item(96, 11).Value = 2532.75  # ALC!K96
$ws.Cells.Item(96, 12).Value = 2275.2  # ALC!L96
$ws.Cells.Item(96, 13).Value = -1159.75  # ALC!M96
$ws.Cells.Item(96, 14).Value = -5021.2  # ALC!N96

$ws.Cells.Item(100, 8).Value = 5575.4287  # ALC!H100
$ws.Cells.Item(100, 9).Value = 1178.2  # ALC!I100
$ws.Cells.Item(100, 10).Value = 8018.3335  # ALC!J100
$ws.Cells.Item(100, 11).Value = 1178.2  # ALC!K100
$ws.Cells.Item(100, 12).Value = 8018.3335  # ALC!L100
$ws.Cells.Item(100, 13).Value = -637.2  # ALC!M100
$ws.Cells.Item(100, 14).Value = -9100.333500000001  # ALC!N100

$ws.Cells.Item(111, 8).Value = 2045.1904  # ALC!H111
$ws.Cells.Item(111, 9).Value = 1619.6923  # ALC!I111
$ws.Cells.Item(111, 11).Value = 4859.0769  # ALC!K111
$ws.Cells.Item(111, 13).Value = -1792.0769  # ALC!M111

$ws.Cells.Item(137, 8).Value = 771409.3  # ALC!H137
$ws.Cells.Item(137, 10).Value = 2582.1  # ALC!J137
$ws.Cells.Item(137, 12).Value = 7746.299999999999  # ALC!L137
$ws.Cells.Item(137, 14).Value = -12846.3  # ALC!N137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2758.2856  # ARM!H2
$ws.Cells.Item(2, 9).Value = 2382.6667  # ARM!I2
$ws.Cells.Item(2, 10).Value = 5012  # ARM!J2
$ws.Cells.Item(2, 11).Value = 2382.6667  # ARM!K2
$ws.Cells.Item(2, 12).Value = 5012  # ARM!L2
$ws.Cells.Item(2, 13).Value = -2269.6667  # ARM!M2
$ws.Cells.Item(2, 14).Value = -5238  # ARM!N2

$ws.Cells.Item(32, 8).Value = 3211.2708  # ARM!H32
$ws.Cells.Item(32, 9).Value = 2871.9512  # ARM!I32
$ws.Cells.Item(32, 11).Value = 2871.9512  # ARM!K32
$ws.Cells.Item(32, 13).Value = -2584.9512  # ARM!M32

$ws.Cells.Item(74, 8).Value = 2765.3704  # ARM!H74
$ws.Cells.Item(74, 10).Value = 3401.0588  # ARM!J74
$ws.Cells.Item(74, 12).Value = 3401.0588  # ARM!L74
$ws.Cells.Item(74, 14).Value = -5149.0588  # ARM!N74

$ws.Cells.Item(77, 8).Value = 2765.3704  # ARM!H77
$ws.Cells.Item(77, 10).Value = 3401.0588  # ARM!J77
$ws.Cells.Item(77, 12).Value = 17005.294  # ARM!L77
$ws.Cells.Item(77, 14).Value = -25741.294  # ARM!N77

$ws.Cells.Item(88, 8).Value = 13378.706  # ARM!H88
$ws.Cells.Item(88, 10).Value = 15022  # ARM!J88
$ws.Cells.Item(88, 12).Value = 15022  # ARM!L88
$ws.Cells.Item(88, 14).Value = -15834  # ARM!N88

$ws.Cells.Item(91, 8).Value = 13378.706  # ARM!H91
$ws.Cells.Item(91, 10).Value = 15022  # ARM!J91
$ws.Cells.Item(91, 12).Value = 15022  # ARM!L91
$ws.Cells.Item(91, 14).Value = -17830  # ARM!N91

$ws.Cells.Item(97, 8).Value = 4789.7036  # ARM!H97
$ws.Cells.Item(97, 9).Value = 5198.773  # ARM!I97
$ws.Cells.Item(97, 11).Value = 5198.773  # ARM!K97
$ws.Cells.Item(97, 13).Value = -4702.773  # ARM!M97

$ws.Cells.Item(116, 8).Value = 2758.2856  # ARM!H116
$ws.Cells.Item(116, 9).Value = 2382.6667  # ARM!I116
$ws.Cells.Item(116, 10).Value = 5012  # ARM!J116
$ws.Cells.Item(116, 11).Value = 2382.6667  # ARM!K116
$ws.Cells.Item(116, 12).Value = 5012  # ARM!L116
$ws.Cells.Item(116, 13).Value = -88.66670000000022  # ARM!M116
$ws.Cells.Item(116, 14).Value = -9600  # ARM!N116

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2758.2856  # BSM!H3
$ws.Cells.Item(3, 9).Value = 2382.6667  # BSM!I3
$ws.Cells.Item(3, 10).Value = 5012  # BSM!J3
$ws.Cells.Item(3, 11).Value = 2382.6667  # BSM!K3
$ws.Cells.Item(3, 12).Value = 5012  # BSM!L3
$ws.Cells.Item(3, 13).Value = -2268.6667  # BSM!M3
$ws.Cells.Item(3, 14).Value = -5240  # BSM!N3

$ws.Cells.Item(26, 8).Value = 36186  # BSM!H26
$ws.Cells.Item(26, 9).Value = 36186  # BSM!I26
$ws.Cells.Item(26, 11).Value = 36186  # BSM!K26
$ws.Cells.Item(26, 13).Value = -35894  # BSM!M26

$ws.Cells.Item(86, 8).Value = 13376.917  # BSM!H86
$ws.Cells.Item(86, 9).Value = 5516.4  # BSM!I86
$ws.Cells.Item(86, 11).Value = 5516.4  # BSM!K86
$ws.Cells.Item(86, 13).Value = -4393.4  # BSM!M86

$ws.Cells.Item(89, 8).Value = 13376.917  # BSM!H89
$ws.Cells.Item(89, 9).Value = 5516.4  # BSM!I89
$ws.Cells.Item(89, 11).Value = 27582  # BSM!K89
$ws.Cells.Item(89, 13).Value = -21966  # BSM!M89

$ws.Cells.Item(94, 8).Value = 39906.332  # BSM!H94
$ws.Cells.Item(94, 9).Value = 9709  # BSM!I94
$ws.Cells.Item(94, 11).Value = 9709  # BSM!K94
$ws.Cells.Item(94, 13).Value = -9258  # BSM!M94

$ws.Cells.Item(99, 8).Value = 4330.3335  # BSM!H99
$ws.Cells.Item(99, 9).Value = 4160.222  # BSM!I99
$ws.Cells.Item(99, 11).Value = 4160.222  # BSM!K99
$ws.Cells.Item(99, 13).Value = -2662.222  # BSM!M99

$ws.Cells.Item(135, 8).Value = 44988.25  # BSM!H135
$ws.Cells.Item(135, 10).Value = 44988.25  # BSM!J135
$ws.Cells.Item(135, 12).Value = 44988.25  # BSM!L135
$ws.Cells.Item(135, 14).Value = -55128.25  # BSM!N135

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4260.3335  # CRP!H16
$ws.Cells.Item(16, 9).Value = 2893  # CRP!I16
$ws.Cells.Item(16, 11).Value = 2893  # CRP!K16
$ws.Cells.Item(16, 13).Value = -2606  # CRP!M16

$ws.Cells.Item(31, 8).Value = 2914.1614  # CRP!H31
$ws.Cells.Item(31, 9).Value = 1947.5  # CRP!I31
$ws.Cells.Item(31, 11).Value = 1947.5  # CRP!K31
$ws.Cells.Item(31, 13).Value = -1652.5  # CRP!M31

$ws.Cells.Item(34, 8).Value = 2914.1614  # CRP!H34
$ws.Cells.Item(34, 9).Value = 1947.5  # CRP!I34
$ws.Cells.Item(34, 11).Value = 1947.5  # CRP!K34
$ws.Cells.Item(34, 13).Value = -1745.5  # CRP!M34

$ws.Cells.Item(50, 8).Value = 54999  # CRP!H50
$ws.Cells.Item(50, 10).Value = 54999  # CRP!J50
$ws.Cells.Item(50, 12).Value = 54999  # CRP!L50
$ws.Cells.Item(50, 14).Value = -56249  # CRP!N50

$ws.Cells.Item(99, 8).Value = 2979.875  # CRP!H99
$ws.Cells.Item(99, 10).Value = 3945  # CRP!J99
$ws.Cells.Item(99, 12).Value = 3945  # CRP!L99
$ws.Cells.Item(99, 14).Value = -6941  # CRP!N99

$ws.Cells.Item(113, 8).Value = 4260.3335  # CRP!H113
$ws.Cells.Item(113, 9).Value = 2893  # CRP!I113
$ws.Cells.Item(113, 11).Value = 2893  # CRP!K113
$ws.Cells.Item(113, 13).Value = -723  # CRP!M113

$ws.Cells.Item(126, 8).Value = 2979.875  # CRP!H126
$ws.Cells.Item(126, 10).Value = 3945  # CRP!J126
$ws.Cells.Item(126, 12).Value = 11835  # CRP!L126
$ws.Cells.Item(126, 14).Value = -16775  # CRP!N126

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 9999  # CUL!H2
$ws.Cells.Item(2, 9).Value = 9999  # CUL!I2
$ws.Cells.Item(2, 10).Value = 0  # CUL!J2
$ws.Cells.Item(2, 11).Value = 59994  # CUL!K2
$ws.Cells.Item(2, 12).ClearContents()  # CUL!L2 remove
$ws.Cells.Item(2, 13).Value = -59881  # CUL!M2
$ws.Cells.Item(2, 14).Value = 0  # CUL!N2

$ws.Cells.Item(5, 8).Value = 1275.3572  # CUL!H5
$ws.Cells.Item(5, 9).Value = 388.72726  # CUL!I5
$ws.Cells.Item(5, 10).Value = 4526.3335  # CUL!J5
$ws.Cells.Item(5, 11).Value = 1166.18178  # CUL!K5
$ws.Cells.Item(5, 12).Value = 13579.0005  # CUL!L5
$ws.Cells.Item(5, 13).Value = -1054.18178  # CUL!M5
$ws.Cells.Item(5, 14).Value = -13803.0005  # CUL!N5

$ws.Cells.Item(107, 8).Value = 616.13336  # CUL!H107
$ws.Cells.Item(107, 9).Value = 184.44444  # CUL!I107
$ws.Cells.Item(107, 10).Value = 1263.6666  # CUL!J107
$ws.Cells.Item(107, 11).Value = 553.33332  # CUL!K107
$ws.Cells.Item(107, 12).Value = 3790.9998  # CUL!L107
$ws.Cells.Item(107, 13).Value = 1366.66668  # CUL!M107
$ws.Cells.Item(107, 14).Value = -7630.9998  # CUL!N107

$ws.Cells.Item(135, 8).Value = 1275.3572  # CUL!H135
$ws.Cells.Item(135, 9).Value = 388.72726  # CUL!I135
$ws.Cells.Item(135, 10).Value = 4526.3335  # CUL!J135
$ws.Cells.Item(135, 11).Value = 3498.54534  # CUL!K135
$ws.Cells.Item(135, 12).Value = 40737.0015  # CUL!L135
$ws.Cells.Item(135, 13).Value = -963.5453400000001  # CUL!M135
$ws.Cells.Item(135, 14).Value = -45807.0015  # CUL!N135

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(31, 8).Value = 18829  # LTW!H31
$ws.Cells.Item(31, 10).Value = 12053  # LTW!J31
$ws.Cells.Item(31, 12).Value = 12053  # LTW!L31
$ws.Cells.Item(31, 14).Value = -12549  # LTW!N31

$ws.Cells.Item(55, 8).Value = 417.375  # LTW!H55
$ws.Cells.Item(55, 9).Value = 256  # LTW!I55
$ws.Cells.Item(55, 10).Value = 772.4  # LTW!J55
$ws.Cells.Item(55, 11).Value = 256  # LTW!K55
$ws.Cells.Item(55, 12).Value = 772.4  # LTW!L55
$ws.Cells.Item(55, 13).Value = -83  # LTW!M55
$ws.Cells.Item(55, 14).Value = -1118.4  # LTW!N55

$ws.Cells.Item(68, 8).Value = 2050.4583  # LTW!H68
$ws.Cells.Item(68, 10).Value = 1999.2  # LTW!J68
$ws.Cells.Item(68, 12).Value = 1999.2  # LTW!L68
$ws.Cells.Item(68, 14).Value = -3497.2  # LTW!N68

$ws.Cells.Item(71, 8).Value = 2050.4583  # LTW!H71
$ws.Cells.Item(71, 10).Value = 1999.2  # LTW!J71
$ws.Cells.Item(71, 12).Value = 9996  # LTW!L71
$ws.Cells.Item(71, 14).Value = -17484  # LTW!N71

$ws.Cells.Item(93, 8).Value = 1666.3334  # LTW!H93
$ws.Cells.Item(93, 9).Value = 1000  # LTW!I93
$ws.Cells.Item(93, 10).Value = 1999.5  # LTW!J93
$ws.Cells.Item(93, 11).Value = 1000  # LTW!K93
$ws.Cells.Item(93, 12).Value = 1999.5  # LTW!L93
$ws.Cells.Item(93, 13).Value = 248  # LTW!M93
$ws.Cells.Item(93, 14).Value = -4495.5  # LTW!N93

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2666.3333  # WVR!H96
$ws.Cells.Item(96, 9).Value = 1999.75  # WVR!I96
$ws.Cells.Item(96, 10).Value = 3999.5  # WVR!J96
$ws.Cells.Item(96, 11).Value = 1999.75  # WVR!K96
$ws.Cells.Item(96, 12).Value = 3999.5  # WVR!L96
$ws.Cells.Item(96, 13).Value = -626.75  # WVR!M96
$ws.Cells.Item(96, 14).Value = -6745.5  # WVR!N96

$ws.Cells.Item(100, 8).Value = 400.64285  # WVR!H100
$ws.Cells.Item(100, 9).Value = 334.0909  # WVR!I100
$ws.Cells.Item(100, 11).Value = 668.1818  # WVR!K100
$ws.Cells.Item(100, 13).Value = -127.1818  # WVR!M100

$ws.Cells.Item(112, 8).Value = 50387  # WVR!H112
$ws.Cells.Item(112, 10).Value = 50387  # WVR!J112
$ws.Cells.Item(112, 12).Value = 50387  # WVR!L112
$ws.Cells.Item(112, 14).Value = -53341  # WVR!N112

$ws.Cells.Item(126, 8).Value = 2321  # WVR!H126
$ws.Cells.Item(126, 9).Value = 2324.5  # WVR!I126
$ws.Cells.Item(126, 10).Value = 2300  # WVR!J126
$ws.Cells.Item(126, 11).Value = 6973.5  # WVR!K126
$ws.Cells.Item(126, 12).Value = 6900  # WVR!L126
$ws.Cells.Item(126, 13).Value = -4503.5  # WVR!M126
$ws.Cells.Item(126, 14).Value = -11840  # WVR!N126
